# Updated symbol list with GitHub Actions — refresh quoted coin prices (col D)
# and a couple of "Worst in 24h" label moves (col E) on Sheet1.
#
# Price cells are stored as text (e.g. "242.65"), matching the source feed's
# formatting (trailing zeros, variable decimal places). A leading apostrophe
# forces Excel to keep them as text instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wb.Worksheets.Item(1).Activate()

# Helper: write a numeric-looking quote as genuine text (leading apostrophe
# forces text entry like a user typing it in), then strip the resulting
# "quote prefix" cell style back to Normal so formatting is untouched.
function Set-TextPrice($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextPrice "D2"  "242.66"
Set-TextPrice "D3"  "23.05"
Set-TextPrice "D4"  "5.425"
Set-TextPrice "D5"  "0.05897"
Set-TextPrice "D6"  "3.437"
Set-TextPrice "D7"  "6.517"
Set-TextPrice "D8"  "0.8101"
Set-TextPrice "D9"  "0.9687"
Set-TextPrice "D10" "0.1419"
Set-TextPrice "D11" "0.07413"
Set-TextPrice "D12" "0.03264"
Set-TextPrice "D13" "0.03054"
Set-TextPrice "D14" "0.09342"
Set-TextPrice "D15" "3.858"
Set-TextPrice "D16" "0.001593"
Set-TextPrice "D17" "0.04674"
Set-TextPrice "D18" "0.0005912"
Set-TextPrice "D19" "0.005857"
Set-TextPrice "D20" "0.001259"
Set-TextPrice "D21" "0.004898"
Set-TextPrice "D22" "0.00006800"
Set-TextPrice "D23" "3.593"
Set-TextPrice "D24" "2.168"
Set-TextPrice "D25" "0.3220"
Set-TextPrice "D27" "0.0002285"
Set-TextPrice "D40" "0.03936"
Set-TextPrice "D41" "0.006187"
Set-TextPrice "D42" "0.1073"
Set-TextPrice "D43" "0.002560"
Set-TextPrice "D44" "0.009140"
Set-TextPrice "D45" "0.00005179"
Set-TextPrice "D47" "0.6662"

# "Worst in 24h" badge moved off row 18 (One) and onto row 47 (CoinbaseStockToken)
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
